$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.693.87'
$ws.Range("E2").Value = '  -1.63%  '

$ws.Range("D3").Value = '3.029.19'
$ws.Range("E3").Value = '  -2.01%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '581.94'
$ws.Range("E5").Value = '  -1.82%  '

$ws.Range("D6").Value = '148.85'
$ws.Range("E6").Value = '  -5.20%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").Value = '0.525'
$ws.Range("E8").Value = '  -3.35%  '

$ws.Range("D9").Value = '3.028.97'
$ws.Range("E9").Value = '  -1.96%  '

$ws.Range("E10").Value = '  -4.10%  '

$ws.Range("D11").Value = '5.65'
$ws.Range("E11").Value = '  -3.27%  '

$ws.Range("E12").Value = '  -2.66%  '

$ws.Range("E13").Value = '  -4.32%  '

$ws.Range("D14").Value = '35.22'
$ws.Range("E14").Value = '  -6.19%  '

$ws.Range("E15").Value = '  +1.73%  '

$ws.Range("D16").Value = '3.530.68'
$ws.Range("E16").Value = '  -2.02%  '

$ws.Range("D17").Value = '62.653.21'
$ws.Range("E17").Value = '  -1.73%  '

$ws.Range("D18").Value = '7.04'
$ws.Range("E18").Value = '  -1.80%  '

$ws.Range("D19").Value = '3.028.26'
$ws.Range("E19").Value = '  -2.04%  '

$ws.Range("D20").Value = '468.22'
$ws.Range("E20").Value = '  -2.52%  '

$ws.Range("D21").Value = '14.03'
$ws.Range("E21").Value = '  -3.93%  '

$ws.Range("D22").Value = '0.691'
$ws.Range("E22").Value = '  -3.08%  '

$ws.Range("E23").Value = '  -2.60%  '

$ws.Range("D24").Value = '2.37'
$ws.Range("E24").Value = '  -1.99%  '

$ws.Range("D25").Value = '81.03'
$ws.Range("E25").Value = '  -0.44%  '

$ws.Range("D26").Value = '12.43'
$ws.Range("E26").Value = '  -3.45%  '

$ws.Range("D27").Value = '10.43'
$ws.Range("E27").Value = '  +1.83%  '

$ws.Range("E28").Value = '  -0.07%  '

$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("D30").Value = '7.20'
$ws.Range("E30").Value = '  -3.78%  '

$ws.Range("E31").Value = '  -2.31%  '

$ws.Range("E32").Value = '  -1.64%  '

$ws.Range("D33").Value = '27.50'
$ws.Range("E33").Value = '  +0.59%  '

$ws.Range("E34").Value = '  -5.61%  '

$ws.Range("E35").Value = '  -1.41%  '

$ws.Range("E36").Value = '  -6.97%  '

$ws.Range("D37").Value = '5.77'
$ws.Range("E37").Value = '  -4.99%  '

$ws.Range("E38").Value = '  -3.30%  '

$ws.Range("D39").Value = '50.28'
$ws.Range("E39").Value = '  -1.33%  '

$ws.Range("E40").Value = '  -15.78%  '

$ws.Range("D41").Value = '8.98'
$ws.Range("E41").Value = '  -3.79%  '

$ws.Range("D42").Value = '421.05'
$ws.Range("E42").Value = '  -5.99%  '

$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D43").Value = '0.280'
$ws.Range("E43").Value = '  -2.94%  '

$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D44").Value = '0.113'
$ws.Range("E44").Value = '  +0.77%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.791.06'
$ws.Range("E45").Value = '  -1.50%  '

$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '0.0355'
$ws.Range("E46").Value = '  -2.17%  '

$ws.Range("D47").Value = '38.00'
$ws.Range("E47").Value = '  -9.67%  '

$ws.Range("D48").Value = '129.56'
$ws.Range("E48").Value = '  -0.76%  '

$ws.Range("E49").Value = '  +0.03%  '

$ws.Range("D50").Value = '24.52'
$ws.Range("E50").Value = '  -4.54%  '

$ws.Range("E51").Value = '  -1.73%  '
